# [base] - [`outputToCloud(resource)`]: support the transferring of output
# artifact to the cloud.
#
# The "#system" sheet stores one lookup table per macro category: row 1
# holds the category name (used as the header / key) and the column below
# it holds the ordered list of command signatures for that category. Named
# ranges (e.g. "base", "web", ...) point at those columns and the
# MacroLibrary sheet's data-validation lists resolve through them
# (target -> category name -> INDIRECT(category) -> command list).
#
# This change:
#   1. Adds a brand new category "text" (inserted alphabetically into the
#      "target" list, right before "web").
#   2. Gives that new category its own column (Y), pushing the existing
#      Y..AD columns (web, webalert, webcookie, ws, ws.async, xml) one
#      column to the right (Z..AE) to make room.
#   3. Adds the new category's single command, spellCheck(var,profile,text).
#   4. Adds a new "base" command, outputToCloud(resource), inserted
#      alphabetically into the existing base list (column E).
#   5. Updates every affected defined name so it keeps pointing at the
#      right range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1. Make room for the new "text" category column by shifting the
#    web / webalert / webcookie / ws / ws.async / xml columns one slot to
#    the right: AD->AE, AC->AD, AB->AC, AA->AB, Z->AA, Y->Z. Walk from the
#    rightmost column to the left so we never overwrite data we still need
#    to read.
# ---------------------------------------------------------------------
$lastRow = 129

$ws.Range("AE1:AE$lastRow").Value = $ws.Range("AD1:AD$lastRow").Value()
$ws.Range("AD1:AD$lastRow").Value = $ws.Range("AC1:AC$lastRow").Value()
$ws.Range("AC1:AC$lastRow").Value = $ws.Range("AB1:AB$lastRow").Value()
$ws.Range("AB1:AB$lastRow").Value = $ws.Range("AA1:AA$lastRow").Value()
$ws.Range("AA1:AA$lastRow").Value = $ws.Range("Z1:Z$lastRow").Value()
$ws.Range("Z1:Z$lastRow").Value  = $ws.Range("Y1:Y$lastRow").Value()

# Column Y is now free (it is a copy of the old "web" column) - clear it
# and write the brand-new "text" category into it: a header plus its one
# command.
$ws.Range("Y1:Y$lastRow").ClearContents()
$ws.Range("Y1").Value = "text"
$ws.Range("Y2").Value = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------
# 2. Insert the new "text" category name into the "target" list (column
#    A), in alphabetical order right before "web" (row 25), by shifting
#    everything from "web" down one row (bottom-up, so nothing is lost).
# ---------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $ws.Range("A" + ($r + 1)).Value = $ws.Range("A$r").Value()
}
$ws.Range("A25").Value = "text"

# ---------------------------------------------------------------------
# 3. Insert the new "outputToCloud(resource)" command into the "base"
#    list (column E), in alphabetical order right before
#    "prependText(var,prependWith)" (row 21), shifting the rest of the
#    list down one row (bottom-up).
# ---------------------------------------------------------------------
for ($r = 38; $r -ge 21; $r--) {
    $ws.Range("E" + ($r + 1)).Value = $ws.Range("E$r").Value()
}
$ws.Range("E21").Value = "outputToCloud(resource)"

# ---------------------------------------------------------------------
# 4. Fix up every defined name so it refers to the right range now that
#    the columns/rows have moved, and add the new "text" name.
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo      = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo    = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo       = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo  = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo        = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo  = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo       = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
